$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell far outside the used range, guaranteed to carry the default (unstyled) format.
$donor = $ws.Cells.Item(1000, 26)

function Set-TextValue($cell, $val) {
    # Force the cell to Text format so a numeric-looking string is not auto-converted
    # to a floating-point number, then restore the original (default) cell format by
    # pasting formats from an always-unstyled donor cell, leaving styles.xml untouched.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $donor.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}

$ws.Cells.Item(2, 4).Value = "68.334.83"
$ws.Cells.Item(2, 5).Value = "  +4.02%  "
$ws.Cells.Item(3, 4).Value = "3.630.94"
$ws.Cells.Item(3, 5).Value = "  +4.52%  "
$ws.Cells.Item(4, 5).Value = "  +0.20%  "
Set-TextValue $ws.Cells.Item(5, 4) "203.27"
$ws.Cells.Item(5, 5).Value = "  +11.76%  "
Set-TextValue $ws.Cells.Item(6, 4) "577.66"
$ws.Cells.Item(6, 5).Value = "  +3.20%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.621"
$ws.Cells.Item(7, 5).Value = "  +3.44%  "
$ws.Cells.Item(8, 5).Value = "  +0.09%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.689"
$ws.Cells.Item(9, 5).Value = "  +6.31%  "
Set-TextValue $ws.Cells.Item(10, 4) "61.31"
$ws.Cells.Item(10, 5).Value = "  +19.39%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.150"
$ws.Cells.Item(11, 5).Value = "  +6.30%  "
Set-TextValue $ws.Cells.Item(12, 4) "0.0000286"
$ws.Cells.Item(12, 5).Value = "  +14.09%  "
Set-TextValue $ws.Cells.Item(13, 4) "10.38"
$ws.Cells.Item(13, 5).Value = "  +9.51%  "
$ws.Cells.Item(14, 4).Value = "4.202.26"
$ws.Cells.Item(14, 5).Value = "  +4.52%  "
$ws.Cells.Item(15, 4).Value = "3.629.04"
$ws.Cells.Item(15, 5).Value = "  +4.74%  "
Set-TextValue $ws.Cells.Item(16, 4) "19.40"
$ws.Cells.Item(16, 5).Value = "  +9.71%  "
$ws.Cells.Item(17, 5).Value = "  +1.71%  "
$ws.Cells.Item(18, 4).Value = "68.152.56"
$ws.Cells.Item(18, 5).Value = "  +4.26%  "
Set-TextValue $ws.Cells.Item(19, 4) "12.41"
$ws.Cells.Item(19, 5).Value = "  +6.47%  "
$ws.Cells.Item(20, 5).Value = "  +4.29%  "
Set-TextValue $ws.Cells.Item(21, 4) "410.77"
$ws.Cells.Item(21, 5).Value = "  +8.81%  "
Set-TextValue $ws.Cells.Item(22, 4) "13.25"
$ws.Cells.Item(22, 5).Value = "  +23.37%  "
$ws.Cells.Item(23, 5).Value = "  +3.86%  "
Set-TextValue $ws.Cells.Item(24, 4) "85.74"
$ws.Cells.Item(24, 5).Value = "  +3.48%  "
$ws.Cells.Item(25, 5).Value = "  +16.13%  "
Set-TextValue $ws.Cells.Item(26, 4) "2.94"
$ws.Cells.Item(26, 5).Value = "  +4.90%  "
Set-TextValue $ws.Cells.Item(27, 4) "12.68"
$ws.Cells.Item(27, 5).Value = "  +7.08%  "
$ws.Cells.Item(28, 5).Value = "  +2.47%  "
Set-TextValue $ws.Cells.Item(29, 4) "9.39"
$ws.Cells.Item(29, 5).Value = "  +9.34%  "
$ws.Cells.Item(30, 5).Value = "  +8.19%  "
Set-TextValue $ws.Cells.Item(31, 4) "31.84"
$ws.Cells.Item(31, 5).Value = "  +5.34%  "
Set-TextValue $ws.Cells.Item(32, 4) "682.79"
$ws.Cells.Item(32, 5).Value = "  +12.19%  "
Set-TextValue $ws.Cells.Item(33, 4) "12.26"
$ws.Cells.Item(33, 5).Value = "  +3.64%  "
$ws.Cells.Item(34, 5).Value = "  +4.78%  "
$ws.Cells.Item(35, 5).Value = "  +2.14%  "
Set-TextValue $ws.Cells.Item(36, 4) "42.06"
$ws.Cells.Item(36, 5).Value = "  +3.66%  "
Set-TextValue $ws.Cells.Item(37, 4) "0.417"
$ws.Cells.Item(37, 5).Value = "  +6.03%  "
$ws.Cells.Item(38, 5).Value = "  -0.39%  "
$ws.Cells.Item(39, 4).Value = "0.0₃0772"
$ws.Cells.Item(39, 5).Value = "  +8.32%  "
Set-TextValue $ws.Cells.Item(40, 4) "3.22"
$ws.Cells.Item(40, 5).Value = "  +18.12%  "
$ws.Cells.Item(41, 5).Value = "  +5.43%  "
$ws.Cells.Item(42, 4).Value = "3.187.27"
$ws.Cells.Item(42, 5).Value = "  +9.79%  "
Set-TextValue $ws.Cells.Item(43, 4) "0.997"
$ws.Cells.Item(43, 5).Value = "  +0.11%  "
Set-TextValue $ws.Cells.Item(44, 4) "2.73"
$ws.Cells.Item(44, 5).Value = "  +12.45%  "
Set-TextValue $ws.Cells.Item(45, 4) "2.87"
$ws.Cells.Item(45, 5).Value = "  +25.99%  "
Set-TextValue $ws.Cells.Item(46, 4) "2.87"
$ws.Cells.Item(46, 5).Value = "  +17.96%  "
Set-TextValue $ws.Cells.Item(47, 4) "0.0417"
$ws.Cells.Item(47, 5).Value = "  +5.92%  "
Set-TextValue $ws.Cells.Item(48, 4) "0.133"
$ws.Cells.Item(48, 5).Value = "  +4.81%  "
$ws.Cells.Item(49, 5).Value = "  +8.23%  "
Set-TextValue $ws.Cells.Item(50, 4) "3.08"
$ws.Cells.Item(50, 5).Value = "  -1.08%  "
Set-TextValue $ws.Cells.Item(51, 4) "139.66"
$ws.Cells.Item(51, 5).Value = "  +1.80%  "

$excel.CutCopyMode = $false
